$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9782881736755371
$ws.Range("B1").Value = 3.052150249481201
$ws.Range("C1").Value = 4.05931282043457
$ws.Range("D1").Value = 2.047739267349243
$ws.Range("E1").Value = 1.217458128929138
